$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.974.74'
$ws.Range("E2").Value = '  +2.60%  '

$ws.Range("D3").Value = '2.691.38'
$ws.Range("E3").Value = '  +1.94%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'522.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("D6").Value = "'148.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.17%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  +1.05%  '

$ws.Range("D9").Value = '2.712.28'
$ws.Range("E9").Value = '  +1.70%  '

$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("D14").Value = '3.164.42'
$ws.Range("E14").Value = '  +2.06%  '

$ws.Range("D15").Value = '60.968.09'
$ws.Range("E15").Value = '  +2.69%  '

$ws.Range("D16").Value = '2.854.79'
$ws.Range("E16").Value = '  +7.24%  '

$ws.Range("D17").Value = "'21.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.09%  '

$ws.Range("E18").Value = '  +0.73%  '

$ws.Range("D19").Value = "'355.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.51%  '

$ws.Range("D20").Value = "'4.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.62%  '

$ws.Range("D21").Value = "'10.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '

$ws.Range("D22").Value = "'6.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.42%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = "'62.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.94%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("E26").Value = '  +4.23%  '

$ws.Range("D27").Value = "'0.989"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.46%  '

$ws.Range("E28").Value = '  +1.06%  '

$ws.Range("D29").Value = "'7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.09%  '

$ws.Range("D30").Value = "'6.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.09%  '

$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'19.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = "'1.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("D34").Value = "'150.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.31%  '

$ws.Range("D35").Value = "'4.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.74%  '

$ws.Range("D36").Value = "'0.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.05%  '

$ws.Range("E37").Value = '  +4.41%  '

$ws.Range("D38").Value = "'1.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.04%  '

$ws.Range("D39").Value = "'0.882"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.46%  '

$ws.Range("D40").Value = "'36.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("D42").Value = "'286.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.97%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = "'0.0995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.36%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'0.615"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'20.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.76%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.144.17'
$ws.Range("E46").Value = '  +7.68%  '

$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = "'0.0543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'4.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.67%  '

$ws.Range("E50").Value = '  +0.60%  '

$ws.Range("D51").Value = "'19.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.54%  '
